$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: assign a value to a cell while forcing it to be stored as text
# (mirrors the original file where numeric-looking values like "0.06" or
# "166108" are stored as strings/inlineStr rather than numbers).
# ---------------------------------------------------------------------------
function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计", before "2022-Q1".
#    This naturally shifts the existing "2022-Q1" and "2020-Q4" sheets/files
#    down, leaving their contents untouched.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTotal)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet.
#    Start by copying the layout/formatting of the "2022-Q1" sheet (same
#    header row + column styling), then overwrite with the new fund data.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q1")
$refSheet.Range("A1:H3").Copy($newSheet.Range("A1"))

# Row 2
$newSheet.Cells.Item(2,1).Value = 0
Set-TextCell $newSheet.Cells.Item(2,2) "166108"
Set-TextCell $newSheet.Cells.Item(2,3) "信澳量化多因子混合（LOF）C"
Set-TextCell $newSheet.Cells.Item(2,4) "0.06"
Set-TextCell $newSheet.Cells.Item(2,5) "28.39"
Set-TextCell $newSheet.Cells.Item(2,6) "0.42"
Set-TextCell $newSheet.Cells.Item(2,7) "0.0003"
$newSheet.Cells.Item(2,8).Value = 3

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
Set-TextCell $newSheet.Cells.Item(3,2) "166107"
Set-TextCell $newSheet.Cells.Item(3,3) "信澳量化多因子混合（LOF）A"
Set-TextCell $newSheet.Cells.Item(3,4) "0.05"
Set-TextCell $newSheet.Cells.Item(3,5) "28.39"
Set-TextCell $newSheet.Cells.Item(3,6) "0.42"
Set-TextCell $newSheet.Cells.Item(3,7) "0.0002"
$newSheet.Cells.Item(3,8).Value = 3

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new row for "2022-Q3"
#    right after the header, pushing the existing "2022-Q1"/"2020-Q4" rows
#    down by one.
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

# Clean up formatting that Excel auto-applies to the inserted row, then
# re-apply the same style the A-column index cells use elsewhere (copied
# from the row that now holds the old "2020-Q4" entry).
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Cells.Item(4,1).Copy($wsTotal.Cells.Item(2,1))

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 2
$wsTotal.Cells.Item(2,4).Value = 0

$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2

Write-Output "2022-Q3 sheet added and 总计 summary updated"
